$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 8): "Tổng xu" column becomes "Tổng doanh thu", and a new
# "Tổng doanh thu xu" column is introduced right after it (branch revenue
# point/coin method support).
$ws.Range("H8").Value = "Tổng doanh thu"
$ws.Range("I8").Value = "Tổng doanh thu xu"

# Move the saved selection (matches the author's last cursor position).
$ws.Range("D14").Select() | Out-Null
